$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.531.11"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "2.673.20"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.99"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.11"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D9").Value = "2.672.80"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.79"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "3.161.96"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "67.468.03"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "2.672.12"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.71"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.68"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.43"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.81"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("E24").Value = "  -4.08%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -4.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("D28").Value = "2.721.76"
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "556.17"
$ws.Range("E31").Value = "  -4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.53"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.57"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("E42").Value = "  -4.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.96"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  -5.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.36"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "0.0₆0299"
$ws.Range("E47").Value = "  -5.69%  "
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.22"
$ws.Range("E49").Value = "  -2.77%  "
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("E51").Value = "  -2.73%  "
